$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1759
$ws.Range("I32").Value = 2500
$ws.Range("J32").Value = 1265
$ws.Range("K32").Value = 2500
$ws.Range("L32").Value = 1265
$ws.Range("M32").Value = -2174
$ws.Range("N32").Value = -1917
$ws.Range("H40").Value = 1018.1177
$ws.Range("I40").Value = 1064.1154
$ws.Range("K40").Value = 1064.1154
$ws.Range("M40").Value = -889.1153999999999
$ws.Range("H51").Value = 2458.889
$ws.Range("J51").Value = 3596
$ws.Range("L51").Value = 3596
$ws.Range("N51").Value = -4564
$ws.Range("H112").Value = 2112.375
$ws.Range("J112").Value = 2483.1667
$ws.Range("L112").Value = 7449.500100000001
$ws.Range("N112").Value = -9665.500100000001
$ws.Range("H125").Value = 1181.7858
$ws.Range("J125").Value = 1568
$ws.Range("L125").Value = 14112
$ws.Range("N125").Value = -19032
$ws.Range("H129").Value = 870.8182
$ws.Range("J129").Value = 893.0769
$ws.Range("L129").Value = 2679.2307
$ws.Range("N129").Value = -12679.2307
$ws.Range("H132").Value = 1237.9592
$ws.Range("I132").Value = 1162.0435
$ws.Range("K132").Value = 3486.1305
$ws.Range("M132").Value = -956.1305000000002
$ws.Range("H137").Value = 1051.2424
$ws.Range("I137").Value = 814.4815
$ws.Range("J137").Value = 2116.6667
$ws.Range("K137").Value = 2443.4445
$ws.Range("L137").Value = 6350.000100000001
$ws.Range("M137").Value = 106.5554999999999
$ws.Range("N137").Value = -11450.0001
$ws.Range("H138").Value = 3525.75
$ws.Range("J138").Value = 4042.6
$ws.Range("L138").Value = 12127.8
$ws.Range("N138").Value = -22407.8
$ws.Range("H141").Value = 1649525.5
$ws.Range("I141").Value = 2546639.8
$ws.Range("K141").Value = 7639919.399999999
$ws.Range("M141").Value = -7634739.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3319.0815
$ws.Range("I32").Value = 2860.2327
$ws.Range("J32").Value = 6607.5
$ws.Range("K32").Value = 2860.2327
$ws.Range("L32").Value = 6607.5
$ws.Range("M32").Value = -2573.2327
$ws.Range("N32").Value = -7181.5
$ws.Range("H61").Value = 2230.878
$ws.Range("I61").Value = 1497.6
$ws.Range("J61").Value = 6508.3335
$ws.Range("K61").Value = 1497.6
$ws.Range("L61").Value = 6508.3335
$ws.Range("M61").Value = -1285.6
$ws.Range("N61").Value = -6932.3335
$ws.Range("H74").Value = 1797.2778
$ws.Range("I74").Value = 1639.4445
$ws.Range("J74").Value = 1955.1111
$ws.Range("K74").Value = 1639.4445
$ws.Range("L74").Value = 1955.1111
$ws.Range("M74").Value = -765.4445000000001
$ws.Range("N74").Value = -3703.1111
$ws.Range("H77").Value = 1797.2778
$ws.Range("I77").Value = 1639.4445
$ws.Range("J77").Value = 1955.1111
$ws.Range("K77").Value = 8197.2225
$ws.Range("L77").Value = 9775.5555
$ws.Range("M77").Value = -3829.2225
$ws.Range("N77").Value = -18511.5555
$ws.Range("H136").Value = 2230.878
$ws.Range("I136").Value = 1497.6
$ws.Range("J136").Value = 6508.3335
$ws.Range("K136").Value = 4492.799999999999
$ws.Range("L136").Value = 19525.0005
$ws.Range("M136").Value = -1942.799999999999
$ws.Range("N136").Value = -24625.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1827.6
$ws.Range("I31").Value = 1666.8572
$ws.Range("K31").Value = 1666.8572
$ws.Range("M31").Value = -1371.8572
$ws.Range("H34").Value = 1827.6
$ws.Range("I34").Value = 1666.8572
$ws.Range("K34").Value = 1666.8572
$ws.Range("M34").Value = -1464.8572
$ws.Range("H58").Value = 1061616.9
$ws.Range("I58").Value = 1450239.5
$ws.Range("K58").Value = 1450239.5
$ws.Range("M58").Value = -1450036.5
$ws.Range("H132").Value = 1499.138
$ws.Range("I132").Value = 1004.6957
$ws.Range("K132").Value = 3014.0871
$ws.Range("M132").Value = -484.0870999999997
$ws.Range("H134").Value = 1962.8529
$ws.Range("I134").Value = 1788.4584
$ws.Range("J134").Value = 2381.4
$ws.Range("K134").Value = 5365.3752
$ws.Range("L134").Value = 7144.200000000001
$ws.Range("M134").Value = -2830.3752
$ws.Range("N134").Value = -12214.2
$ws.Range("H136").Value = 1061616.9
$ws.Range("I136").Value = 1450239.5
$ws.Range("K136").Value = 4350718.5
$ws.Range("M136").Value = -4348168.5
$ws.Range("H139").Value = 10000
$ws.Range("I139").Value = 10000
$ws.Range("K139").Value = 10000
$ws.Range("M139").Value = -4860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 3634.818
$ws.Range("J110").Value = 3795.6
$ws.Range("L110").Value = 11386.8
$ws.Range("N110").Value = -19566.8
$ws.Range("H113").Value = 6029.1055
$ws.Range("I113").Value = 14756.857
$ws.Range("J113").Value = 937.9167
$ws.Range("K113").Value = 44270.571
$ws.Range("L113").Value = 2813.7501
$ws.Range("M113").Value = -42100.571
$ws.Range("N113").Value = -7153.7501
$ws.Range("H131").Value = 12850.717
$ws.Range("J131").Value = 14035.869
$ws.Range("L131").Value = 42107.607
$ws.Range("N131").Value = -52187.607
$ws.Range("H132").Value = 910
$ws.Range("J132").Value = 1082.6666
$ws.Range("L132").Value = 9743.999400000001
$ws.Range("N132").Value = -14803.9994
$ws.Range("H134").Value = 2967.5715
$ws.Range("I134").Value = 2194.5
$ws.Range("J134").Value = 3998.3333
$ws.Range("K134").Value = 6583.5
$ws.Range("L134").Value = 11994.9999
$ws.Range("M134").Value = -1513.5
$ws.Range("N134").Value = -22134.9999
$ws.Range("H139").Value = 13290.223
$ws.Range("I139").Value = 16230.286
$ws.Range("K139").Value = 48690.858
$ws.Range("M139").Value = -43550.858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1508.1
$ws.Range("I122").Value = 1220.25
$ws.Range("J122").Value = 1939.875
$ws.Range("K122").Value = 3660.75
$ws.Range("L122").Value = 5819.625
$ws.Range("M122").Value = -1210.75
$ws.Range("N122").Value = -10719.625
$ws.Range("H134").Value = 28442
$ws.Range("J134").Value = 28442
$ws.Range("L134").Value = 85326
$ws.Range("N134").Value = -90396
$ws.Range("H141").Value = 33248.75
$ws.Range("J141").Value = 33248.75
$ws.Range("L141").Value = 33248.75
$ws.Range("N141").Value = -43608.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4750
$ws.Range("I22").Value = 10000
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 10000
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -9705
$ws.Range("N22").Value = -3590
$ws.Range("H27").Value = 4750
$ws.Range("I27").Value = 10000
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 10000
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -9893
$ws.Range("N27").Value = -3214
$ws.Range("H40").Value = 7388.375
$ws.Range("I40").Value = 8092.375
$ws.Range("J40").Value = 5980.375
$ws.Range("K40").Value = 8092.375
$ws.Range("L40").Value = 5980.375
$ws.Range("M40").Value = -7956.375
$ws.Range("N40").Value = -6252.375
$ws.Range("H82").Value = 1248.3572
$ws.Range("I82").Value = 1052.8
$ws.Range("K82").Value = 1052.8
$ws.Range("M82").Value = -691.8
$ws.Range("H85").Value = 1248.3572
$ws.Range("I85").Value = 1052.8
$ws.Range("K85").Value = 1052.8
$ws.Range("M85").Value = 195.2
$ws.Range("H132").Value = 1383.96
$ws.Range("I132").Value = 965.30554
$ws.Range("K132").Value = 2895.91662
$ws.Range("M132").Value = -365.91662
$ws.Range("H135").Value = 33325.8
$ws.Range("J135").Value = 33325.8
$ws.Range("L135").Value = 33325.8
$ws.Range("N135").Value = -43465.8
$ws.Range("H136").Value = 2141.3333
$ws.Range("I136").Value = 1327.1613
$ws.Range("J136").Value = 5296.25
$ws.Range("K136").Value = 3981.4839
$ws.Range("L136").Value = 15888.75
$ws.Range("M136").Value = -1431.4839
$ws.Range("N136").Value = -20988.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 951.5
$ws.Range("I107").Value = 732.4
$ws.Range("J107").Value = 1316.6666
$ws.Range("K107").Value = 2197.2
$ws.Range("L107").Value = 3949.9998
$ws.Range("M107").Value = -277.1999999999998
$ws.Range("N107").Value = -7789.9998
$ws.Range("H113").Value = 885.4
$ws.Range("I113").Value = 481.75
$ws.Range("K113").Value = 1445.25
$ws.Range("M113").Value = 724.75
$ws.Range("H126").Value = 13487.75
$ws.Range("I126").Value = 18884
$ws.Range("J126").Value = 8091.5
$ws.Range("K126").Value = 56652
$ws.Range("L126").Value = 24274.5
$ws.Range("M126").Value = -54182
$ws.Range("N126").Value = -29214.5
$ws.Range("H132").Value = 1596.8
$ws.Range("I132").Value = 944.3913
$ws.Range("K132").Value = 2833.1739
$ws.Range("M132").Value = -303.1738999999998
$ws.Range("H136").Value = 34725444
$ws.Range("I136").Value = 50507916
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 151523748
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -151521198
$ws.Range("N136").Value = -17100
